$d = $word.ActiveDocument

$pairs = @(
    @("50×44=", "56×71="),
    @("70×64=", "15×46="),
    @("44×33=", "85×50="),
    @("50×20=", "84×66="),
    @("83×30=", "17×40="),
    @("39×95=", "76×89="),
    @("34×13=", "65×98="),
    @("18×18=", "96×15="),
    @("31×18=", "77×53="),
    @("12×90=", "41×69="),
    @("36×36=", "96×66="),
    @("29×56=", "59×21="),
    @("65×96=", "20×50="),
    @("48×15=", "32×84="),
    @("21×86=", "23×53="),
    @("49×33=", "42×68="),
    @("50×94=", "19×26="),
    @("88×12=", "30×71="),
    @("35×60=", "18×15="),
    @("72×26=", "90×32="),
    @("36×50=", "32×98="),
    @("20×16=", "48×63="),
    @("71×64=", "52×35="),
    @("69×93=", "33×16="),
    @("73×83=", "90×16=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
